# Apply the "Updated cryptos list" data refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price strings (e.g. "1.00", "594.06") need a leading
# apostrophe (quote-prefix) so Excel stores them as text, exactly as the
# source data feed renders them (trailing zeros, fixed decimals, etc.)
# rather than auto-coercing them into numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.482.35"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.668.61"
$ws.Range("E3").Value = "  +3.67%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'594.06"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'144.09"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.668.31"
$ws.Range("E9").Value = "  +3.76%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.89%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.81%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.92%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +2.21%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'27.54"
$ws.Range("E14").Value = "  +2.63%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.147.03"
$ws.Range("E15").Value = "  +3.76%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "63.379.05"
$ws.Range("E16").Value = "  +1.07%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.0000144"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.664.02"
$ws.Range("E18").Value = "  +3.40%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'11.56"
$ws.Range("E19").Value = "  +4.49%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'4.43"
$ws.Range("E20").Value = "  +2.29%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'339.33"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.83"
$ws.Range("E22").Value = "  +2.93%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.04%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'67.52"
$ws.Range("E24").Value = "  +1.04%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  +7.26%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  +1.75%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +1.20%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +3.83%  "

# Row 31 - Aptos
$ws.Range("D31").Value = "'7.87"
$ws.Range("E31").Value = "  -0.21%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +12.88%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "'1.98"
$ws.Range("E33").Value = "  +2.94%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0813"
$ws.Range("E34").Value = "  +2.39%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'172.99"
$ws.Range("E35").Value = "  -1.96%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'5.14"
$ws.Range("E36").Value = "  +15.93%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  +2.56%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.13%  "

# Row 39 - EthereumClassic
$ws.Range("D39").Value = "'19.20"
$ws.Range("E39").Value = "  +1.99%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +8.89%  "

# Row 41 - Aave
$ws.Range("D41").Value = "'174.91"
$ws.Range("E41").Value = "  +11.35%  "

# Row 42 - USDe
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'40.19"
$ws.Range("E43").Value = "  +0.17%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'3.77"
$ws.Range("E44").Value = "  +2.13%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'22.34"
$ws.Range("E45").Value = "  +6.12%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0564"
$ws.Range("E46").Value = "  +5.86%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +0.80%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +2.97%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +0.71%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'18.80"
$ws.Range("E50").Value = "  +4.85%  "

# Row 51 - dogwifhat
$ws.Range("D51").Value = "'1.73"
$ws.Range("E51").Value = "  +3.66%  "

# Row 29 - was Binance-PegBSC-USD, now Bittensor (rows 29/30 swapped)
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'541.82"
$ws.Range("E29").Value = "  +19.37%  "

# Row 30 - was Bittensor, now Binance-PegBSC-USD
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.06%  "

